$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 16 with a "proceed" / "PROCEED" key-value pair, following
# the same pattern as the other rows in the sheet.
$ws.Range("A16").Value = "proceed"
$ws.Range("B16").Value = "PROCEED"

# Match formatting used by the other value cells in column B (wrap text style).
$ws.Range("B16").WrapText = $true

# Update the selected cell to reflect the new active cell after the edit.
$ws.Range("A17").Select()
